# Corrected excel sheets for application fix issues
$wb = $excel.ActiveWorkbook

$wsModify  = $wb.Worksheets.Item("Modify Transaction")
$wsSummary = $wb.Worksheets.Item("Summary")

# Summary sheet: correct the "Original" and "Written Off" amounts on row 3
$wsSummary.Range("A3").Value = 297.55
$wsSummary.Range("D3").Value = 297.55

# Update selections on both sheets, then make "Summary" the active tab
$wsModify.Activate()
$wsModify.Range("B3").Select()

$wsSummary.Activate()
$wsSummary.Range("D4").Select()
